$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.897.96'
$ws.Range("E2").Value = '  +4.96%  '
$ws.Range("D3").Value = '2.232.34'
$ws.Range("E3").Value = '  +1.95%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.37'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.94%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.405'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0911'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.31%  '
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '2.568.28'
$ws.Range("E13").Value = '  +1.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.805'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("E17").Value = '  +1.46%  '
$ws.Range("D18").Value = '2.241.92'
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("D19").Value = '41.789.29'
$ws.Range("E19").Value = '  +4.72%  '
$ws.Range("D20").Value = '0.0₃0909'
$ws.Range("E20").Value = '  +0.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.22%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.144'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.47%  '
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.06'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.62%  '
$ws.Range("E35").Value = '  +2.85%  '
$ws.Range("E36").Value = '  +2.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.65'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.74'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.36'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000259'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +33.12%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("E42").Value = '  +5.61%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.69%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.39%  '
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '99.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0959'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.47%  '
$ws.Range("D48").Value = '1.482.74'
$ws.Range("E48").Value = '  -2.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.66%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.76%  '
